$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-17 22:51:59"
$wsZh.Range("H2").Value = "2016-03-17 22:52:17"

# de-de sheet: update Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-17 22:52:03"
$wsDe.Range("H2").Value = "2016-03-17 22:52:23"
